# Updates cryptocurrency price (D) and volume-change (E) columns
# to reflect the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.537.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.128.28"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("E9").Value = "  +2.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.45%  "

$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").Value = "  +2.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.675.48"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.17%  "

$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.704.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.134.75"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.90"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "354.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.507"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.81%  "

$ws.Range("E26").Value = "  +1.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0925"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("E29").Value = "  +2.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.32"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.19%  "

$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.90"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.37"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.16"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.96"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.49%  "

$ws.Range("E39").Value = "  -3.26%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("E41").Value = "  +6.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.702"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.173.86"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("E45").Value = "  +4.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.59"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.326.46"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.979"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.33"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.12%  "

